# "Cambios en el Documento" - desarrollo de la metodologia del documento.
#
# 1) Hoja1: se eliminan las columnas auxiliares F:Q (tabla dinamica de
#    apoyo con UNIQUE/TRANSPOSE/SEQUENCE/LARGE que ya no se usa).
# 2) no2: se reordena (ordena) la tabla por la columna "Valor" (D) de
#    mayor a menor, pues la ultima fila se habia agregado sin ordenar.
# 3) Se ajustan las vistas/selecciones y la hoja activa (pm10 pasa a ser
#    la hoja activa en lugar de no2).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("pm10")
$ws3 = $wb.Worksheets.Item("pm25")
$ws4 = $wb.Worksheets.Item("no2")

# --- Hoja1: borrar la tabla auxiliar de columnas F:Q ------------------
# Encabezados de la tabla transpuesta (G1:Q1, formula de array
# TRANSPOSE(UNIQUE(...))) y la columna indice (F) + columna de valores
# (G) usada por LARGE(...).
$ws1.Range("G1:Q1").ClearContents()
$ws1.Range("F1:F14").ClearContents()
$ws1.Range("G2:G12").ClearContents()

# --- no2: reordenar la tabla por la columna Valor (D), descendente ---
$rango = $ws4.Range("A1:D14")
$clave = $ws4.Range("D1:D14")
[void]$rango.Sort($clave, 2)

# --- Vistas / selecciones ---------------------------------------------
# Hoja1 vuelve a la seleccion por defecto (A1).
[void]$ws1.Range("A1").Select()

# pm10 pasa a ser la hoja activa, con la celda E20 seleccionada.
[void]$ws2.Activate()
[void]$ws2.Range("E20").Select()

# pm25 queda con el rango de la tabla completo seleccionado (pero sin
# ser la hoja activa).
[void]$ws3.Range("A1:D14").Select()

# no2 deja de ser la hoja activa y su seleccion queda en D14.
[void]$ws4.Range("D14").Select()

# pm10 es la hoja que debe quedar activa/visible al guardar.
[void]$ws2.Activate()
